# Refresh the cryptocurrency price / volume(1h) snapshot and restore the
# two row pairs whose coins had been swapped (rows 20/21 and 37/38).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.899.62'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.28%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.632.15'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.50%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.95'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.71%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5113'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.49%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2571'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.44%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06339'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.45%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.44'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.19%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07774'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.13%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.272'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.29%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.638.53'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.30%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.857.95'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.45%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5506'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.19%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.84'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.64%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7640'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.24%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.932.57'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.11%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.05%  '
# Row 20
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.51'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.80%  '
# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.414'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.33%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.854'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.68%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.031'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.13%  '
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.891'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.33%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.76'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.64%  '
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.01%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.753'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.35%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.54'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.05%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.241'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.58%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04875'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.50%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.235'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.39%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.183'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.34%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.540'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.88%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.373'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.35%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8972'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.45%  '
# Row 37
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.540'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.52%  '
# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5504'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.36%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.116.99'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.40%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01559'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.29%  '
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.000'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.15%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.581'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.65%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7965'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.93%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.55'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -1.41%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.767.56'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.54%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈114'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -10.12%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4445'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.84%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.36%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.70'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.41%  '
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.55%  '
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.527'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.72%  '
